$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18, shifting existing rows 18-98 down to 19-99.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new weekly record.
$ws.Cells.Item(18, 1).Value  = 6
$ws.Cells.Item(18, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(18, 3).Value  = "Metropolitana"
$ws.Cells.Item(18, 4).Value  = 44802
$ws.Cells.Item(18, 5).Value  = 13
$ws.Cells.Item(18, 6).Value  = 100114007
$ws.Cells.Item(18, 7).Value  = "Jengibre"
$ws.Cells.Item(18, 8).Value  = "Sin especificar"
$ws.Cells.Item(18, 9).Value  = "Primera"
$ws.Cells.Item(18, 10).Value = 500
$ws.Cells.Item(18, 11).Value = 13000
$ws.Cells.Item(18, 12).Value = 15000
$ws.Cells.Item(18, 13).Value = 14080
$ws.Cells.Item(18, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(18, 15).Value = "Perú"
$ws.Cells.Item(18, 16).Value = 1083
$ws.Cells.Item(18, 17).Value = 13
$ws.Cells.Item(18, 18).Value = "Hortaliza"

# Match the date number format used by the rest of column D.
$ws.Cells.Item(18, 4).NumberFormat = $ws.Cells.Item(19, 4).NumberFormat
